$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate the formatting of A2 (fill/border + quote-prefixed text style)
# down onto the new row, then enter the new "120" line item as text so it
# is stored as a shared string just like the existing "100" entry.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A3").Value = "'120"

# Leave the selection where the user would land after typing the entry.
$ws.Range("B3").Select()
